$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.600.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.444.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.34%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '672.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.08%  '
$ws.Range("E7").Value = '  +8.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.465'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +15.91%  '
$ws.Range("E9").Value = '  +21.08%  '
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.444.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.221'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000275'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +13.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.757.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.075.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +33.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.437.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.27%  '
$ws.Range("E20").Value = '  +14.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '538.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +12.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000217'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.439'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +50.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +16.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.618.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.150'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +13.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +19.61%  '
$ws.Range("E32").Value = '  +8.31%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '30.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.565'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +23.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +15.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.61%  '
$ws.Range("E39").Value = '  +9.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '534.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0439'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +35.49%  '
$ws.Range("E44").Value = '  +8.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +15.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +17.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +17.07%  '
